# Updates the cryptos list: price (D) and 1h volume change (E) columns,
# plus a resort of rows 37-42 (coin name/link B,C columns included).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.517.49"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "'2.930.34"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'199.69"
$ws.Range("E5").Value = "  +2.18%  "

$ws.Range("D6").Value = "'593.66"
$ws.Range("E6").Value = "  -1.13%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("D9").Value = "'0.194"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").Value = "'2.928.29"
$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("D11").Value = "'0.446"
$ws.Range("E11").Value = "  +13.48%  "

$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "'3.464.98"
$ws.Range("E14").Value = "  +2.48%  "

$ws.Range("D15").Value = "'76.250.33"
$ws.Range("E15").Value = "  +0.14%  "

$ws.Range("D16").Value = "'28.02"
$ws.Range("E16").Value = "  +1.83%  "

$ws.Range("D17").Value = "'0.0000187"
$ws.Range("E17").Value = "  -1.17%  "

$ws.Range("D18").Value = "'2.910.93"
$ws.Range("E18").Value = "  +1.96%  "

$ws.Range("D19").Value = "'13.22"
$ws.Range("E19").Value = "  +6.23%  "

$ws.Range("D20").Value = "'8.59"
$ws.Range("E20").Value = "  -5.83%  "

$ws.Range("D21").Value = "'369.52"
$ws.Range("E21").Value = "  -3.50%  "

$ws.Range("E22").Value = "  +4.23%  "

$ws.Range("E23").Value = "  -4.15%  "

$ws.Range("D24").Value = "'72.01"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").Value = "'3.069.47"
$ws.Range("E26").Value = "  +2.03%  "

$ws.Range("D27").Value = "'4.21"
$ws.Range("E27").Value = "  -0.21%  "

$ws.Range("D28").Value = "'9.68"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").Value = "'0.0000106"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'8.09"
$ws.Range("E31").Value = "  +4.87%  "

$ws.Range("E32").Value = "  -3.45%  "

$ws.Range("D33").Value = "'492.66"
$ws.Range("E33").Value = "  -4.03%  "

$ws.Range("E34").Value = "  -0.08%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").Value = "'166.87"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("B37").Value = "Cronos"
$ws.Range("C37").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D37").Value = "'0.111"
$ws.Range("E37").Value = "  +21.09%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'20.08"
$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "'0.392"
$ws.Range("E39").Value = "  +13.98%  "

$ws.Range("D40").Value = "'19.72"
$ws.Range("E40").Value = "  +1.23%  "

$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.110"
$ws.Range("E42").Value = "  -7.84%  "

$ws.Range("D43").Value = "'179.56"
$ws.Range("E43").Value = "  -3.47%  "

$ws.Range("D44").Value = "'4.87"
$ws.Range("E44").Value = "  -4.27%  "

$ws.Range("E45").Value = "  -2.92%  "

$ws.Range("D46").Value = "'40.08"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("E47").Value = "  -5.38%  "

$ws.Range("D48").Value = "'0.587"
$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("D49").Value = "'3.86"
$ws.Range("E49").Value = "  +3.12%  "

$ws.Range("D50").Value = "'2.30"
$ws.Range("E50").Value = "  -3.96%  "

$ws.Range("D51").Value = "'22.43"
$ws.Range("E51").Value = "  +4.18%  "
